$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3, 5, 6, 7, 10, 11, 12, 13, 14, 15, 17, 18, 19, 29, 30, 32, 33, 34, 37, 38, 39, 40, 41, 42, 44, 45, 46, 56, 57, 59, 60, 61, 64, 65, 66, 67, 68, 69, 71, 72, 73, 84, 85, 86, 87, 88, 89, 93, 95, 96, 97, 110, 111, 112, 113, 114, 115, 119, 121, 122, 123, 136, 137, 138, 139, 140, 141, 145, 147, 148, 149)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2
    $parts = $v -split ", "
    $n = $parts.Length
    if ($n -gt 1) {
        $rev = $parts[($n-1)..0]
        $joined = $rev -join ", "
        $cell.Value = $joined
    }
}
